# Generate Report for handoff
# Updates the "Latest Handoff Datetime" for the 4a4f9c03-... file (row 4)
# on both the zh-cn and de-de localization-status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-26 11:31:58"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-26 11:32:10"
